# The deck currently carries its custom "Integral" theme (colour scheme)
# on the part backing $p.SlideMaster.Theme. The authored edit swaps the
# two theme parts so the master's theme becomes the stock "Office Theme"
# colour palette. PowerPoint's object model has no call that rewrites a
# theme part wholesale, but it does let automation repaint every slot of
# a master's ThemeColorScheme - which is how a "switch the deck's colours
# back to the Office defaults" edit is actually performed through COM,
# and it recolours every slide/layout object bound to those scheme-colour
# slots (accent1, dk2, ...) at the same time.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Office Theme palette, in ThemeColorScheme slot order (values are the
# packed BGR integers PowerPoint's ColorFormat.RGB setter expects, i.e.
# what VBA's RGB(r,g,b) would return for each hex colour below):
$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
